$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("des_Bildungstand")

$ws.Range("C6").Value = 34.38
$ws.Range("G6").Value = 39.27

$ws.Range("C7").Value = 35.09
$ws.Range("F7").Value = 34.7
$ws.Range("G7").Value = 40.1
$ws.Range("I7").Value = 5.94
